$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values could be misread as numbers by Excel's type
# inference (e.g. "0.999", "103.48"). Force them to remain plain text by
# temporarily setting the cell number format to Text, then restore the
# cell style afterwards so no visible formatting changes are left behind.
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '389.01'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '103.48'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '37.11'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0862'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '18.69'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.81'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.984'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '3.20'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '70.05'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '268.95'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '8.17'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '27.14'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.21'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.109'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '35.82'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '50.49'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.999'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.41'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.291'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.89'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '17.00'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.59'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '127.42'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '22.13'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.48'
$ws.Cells.Item(47, 4).Style = "Normal"

# Remaining cells are unambiguously text already (percentages with padding,
# or price strings using "." as a thousands separator), so a direct .Value
# assignment keeps them as text without any extra bookkeeping.
$ws.Cells.Item(2, 4).Value = '51.735.37'
$ws.Cells.Item(2, 5).Value = '  +0.43%  '
$ws.Cells.Item(3, 4).Value = '3.099.41'
$ws.Cells.Item(3, 5).Value = '  +3.85%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 5).Value = '  +1.82%  '
$ws.Cells.Item(6, 5).Value = '  -0.73%  '
$ws.Cells.Item(7, 5).Value = '  -0.02%  '
$ws.Cells.Item(8, 5).Value = '  -0.02%  '
$ws.Cells.Item(9, 5).Value = '  -0.72%  '
$ws.Cells.Item(10, 5).Value = '  +1.23%  '
$ws.Cells.Item(11, 5).Value = '  +0.09%  '
$ws.Cells.Item(12, 5).Value = '  +0.45%  '
$ws.Cells.Item(13, 4).Value = '3.579.92'
$ws.Cells.Item(13, 5).Value = '  +3.48%  '
$ws.Cells.Item(14, 5).Value = '  +1.27%  '
$ws.Cells.Item(15, 5).Value = '  -0.09%  '
$ws.Cells.Item(16, 4).Value = '3.093.38'
$ws.Cells.Item(16, 5).Value = '  +3.50%  '
$ws.Cells.Item(17, 5).Value = '  -1.11%  '
$ws.Cells.Item(18, 5).Value = '  -4.58%  '
$ws.Cells.Item(19, 4).Value = '51.846.20'
$ws.Cells.Item(19, 5).Value = '  +0.65%  '
$ws.Cells.Item(20, 5).Value = '  +3.41%  '
$ws.Cells.Item(21, 5).Value = '  -0.96%  '
$ws.Cells.Item(22, 4).Value = '0.0₃0967'
$ws.Cells.Item(22, 5).Value = '  +0.48%  '
$ws.Cells.Item(23, 5).Value = '  -0.36%  '
$ws.Cells.Item(24, 5).Value = '  +0.63%  '
$ws.Cells.Item(25, 5).Value = '  -2.81%  '
$ws.Cells.Item(26, 5).Value = '  +3.96%  '
$ws.Cells.Item(27, 5).Value = '  +4.07%  '
$ws.Cells.Item(28, 5).Value = '  +0.75%  '
$ws.Cells.Item(29, 5).Value = '  -1.14%  '
$ws.Cells.Item(30, 5).Value = '  +0.08%  '
$ws.Cells.Item(31, 5).Value = '  -0.11%  '
$ws.Cells.Item(32, 5).Value = '  -0.48%  '
$ws.Cells.Item(33, 5).Value = '  +3.32%  '
$ws.Cells.Item(34, 5).Value = '  +0.70%  '
$ws.Cells.Item(35, 5).Value = '  -1.76%  '
$ws.Cells.Item(36, 5).Value = '  +0.50%  '
$ws.Cells.Item(37, 5).Value = '  -0.15%  '
$ws.Cells.Item(38, 5).Value = '  +3.65%  '
$ws.Cells.Item(39, 5).Value = '  +7.08%  '
$ws.Cells.Item(40, 5).Value = '  +2.20%  '
$ws.Cells.Item(41, 5).Value = '  +0.21%  '
$ws.Cells.Item(42, 5).Value = '  +0.78%  '
$ws.Cells.Item(43, 5).Value = '  -0.33%  '
$ws.Cells.Item(44, 5).Value = '  +1.08%  '
$ws.Cells.Item(45, 5).Value = '  -3.20%  '
$ws.Cells.Item(46, 5).Value = '  +3.54%  '
$ws.Cells.Item(47, 5).Value = '  +5.01%  '
$ws.Cells.Item(48, 5).Value = '  +2.40%  '
$ws.Cells.Item(49, 4).Value = '2.047.36'
$ws.Cells.Item(49, 5).Value = '  +1.17%  '
$ws.Cells.Item(50, 4).Value = '3.403.10'
$ws.Cells.Item(50, 5).Value = '  +3.78%  '
$ws.Cells.Item(51, 5).Value = '  +6.32%  '
